$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it gets the same bold/centered/bordered style (s="1"),
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column data (era data updated).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
